{"js": "// Update the date line and the 25 \"two-digit \u00f7 one-digit\" practice\n// problems in the table to the new day's values.\n\n// 1) Header date paragraph: \"2026-01-20 Tuesday\" -> \"2026-01-21 Wednesday\"\nconst body = context.document.body;\nconst paras = body.paragraphs;\nparas.load(\"items\");\nawait context.sync();\n\nconst firstPara = paras.items[0];\nfirstPara.load(\"text\");\nawait context.sync();\nif (firstPara.text.trim() === \"2026-01-20 Tuesday\") {\n  firstPara.insertText(\"2026-01-21 Wednesday\", Word.InsertLocation.replace);\n}\n\n// 2) Table of division problems. Data lives in rows 0, 4, 8, 12, 16 of the\n// single table (the rows in between are blank rows left for handwritten\n// answers), 5 cells per row.\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Row index (within the table) -> [oldText -> newText] for each of its 5 cells.\nconst rowUpdates = {\n  0: [\"58\u00f72=\", \"90\u00f74=\", \"54\u00f73=\", \"10\u00f79=\", \"58\u00f77=\"],\n  4: [\"21\u00f76=\", \"33\u00f74=\", \"56\u00f75=\", \"20\u00f76=\", \"80\u00f78=\"],\n  8: [\"20\u00f77=\", \"80\u00f76=\", \"88\u00f72=\", \"71\u00f77=\", \"25\u00f77=\"],\n  12: [\"43\u00f79=\", \"68\u00f73=\", \"83\u00f73=\", \"10\u00f77=\", \"28\u00f72=\"],\n  16: [\"98\u00f72=\", \"18\u00f77=\", \"96\u00f79=\", \"59\u00f74=\", \"69\u00f74=\"]\n};\n\nfor (const rowIndexStr of Object.keys(rowUpdates)) {\n  const rowIndex = parseInt(rowIndexStr, 10);\n  const newValues = rowUpdates[rowIndex];\n  for (let cellIndex = 0; cellIndex < newValues.length; cellIndex++) {\n    const cell = table.getCell(rowIndex, cellIndex);\n    const cellParas = cell.body.paragraphs;\n    cellParas.load(\"items\");\n    await context.sync();\n    const p = cellParas.items[0];\n    p.insertText(newValues[cellIndex], Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date line and the 25 \"two-digit \u00f7 one-digit\" practice\n# problems in the table to the new day's values.\n\n$d = $word.ActiveDocument\n\n# 1) Header date paragraph: \"2026-01-20 Tuesday\" -> \"2026-01-21 Wednesday\"\n$d.Paragraphs.Item(1).Range.Text = \"2026-01-21 Wednesday\"\n\n# 2) Table of division problems. Data lives in rows 1, 5, 9, 13, 17\n# (1-based) of the single table (the rows in between are blank rows left\n# for handwritten answers), 5 cells per row.\n$t = $d.Tables.Item(1)\n\n$rowUpdates = @{\n    1  = @(\"58\u00f72=\", \"90\u00f74=\", \"54\u00f73=\", \"10\u00f79=\", \"58\u00f77=\")\n    5  = @(\"21\u00f76=\", \"33\u00f74=\", \"56\u00f75=\", \"20\u00f76=\", \"80\u00f78=\")\n    9  = @(\"20\u00f77=\", \"80\u00f76=\", \"88\u00f72=\", \"71\u00f77=\", \"25\u00f77=\")\n    13 = @(\"43\u00f79=\", \"68\u00f73=\", \"83\u00f73=\", \"10\u00f77=\", \"28\u00f72=\")\n    17 = @(\"98\u00f72=\", \"18\u00f77=\", \"96\u00f79=\", \"59\u00f74=\", \"69\u00f74=\")\n}\n\nforeach ($rowIndex in $rowUpdates.Keys) {\n    $values = $rowUpdates[$rowIndex]\n    for ($col = 1; $col -le $values.Count; $col++) {\n        $t.Cell($rowIndex, $col).Range.Text = $values[$col - 1]\n    }\n}\n"}
